$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        $Sheet,
        [string]$Addr,
        [string]$Val
    )
    $rng = $Sheet.Range($Addr)
    # Force the cell to a text number format first so that assigning a
    # numeric-looking string (e.g. "307.69" or "0.92%") is stored as text
    # instead of being auto-converted into a numeric value by Excel.
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    # Drop the temporary text format again so the cell keeps the default
    # (unstyled) formatting it had before, matching the original workbook.
    $rng.ClearFormats()
}

Set-TextCell $ws "D2" "307.69"
Set-TextCell $ws "E2" "0.92%"

Set-TextCell $ws "D3" "36.24"
Set-TextCell $ws "E3" "1.38%"

Set-TextCell $ws "D4" "5.060"
Set-TextCell $ws "E4" "1.36%"

Set-TextCell $ws "D5" "0.08115"

Set-TextCell $ws "D6" "1.990"
Set-TextCell $ws "E6" "4.71%"

Set-TextCell $ws "D7" "7.858"
Set-TextCell $ws "E7" "-0.31%"

Set-TextCell $ws "D8" "0.9282"

Set-TextCell $ws "D9" "0.1465"
Set-TextCell $ws "E9" "18.86%"

Set-TextCell $ws "D10" "0.1931"
Set-TextCell $ws "E10" "1.62%"

Set-TextCell $ws "D11" "0.09137"
Set-TextCell $ws "E11" "-0.70%"

Set-TextCell $ws "D12" "0.03461"
Set-TextCell $ws "E12" "-1.15%"

Set-TextCell $ws "D13" "0.09892"
Set-TextCell $ws "E13" "-0.22%"

Set-TextCell $ws "D14" "0.001405"
Set-TextCell $ws "E14" "-0.79%"

Set-TextCell $ws "D15" "0.006094"
Set-TextCell $ws "E15" "-2.95%"

Set-TextCell $ws "D16" "3.842"
Set-TextCell $ws "E16" "6.45%"

Set-TextCell $ws "D17" "4.163"
Set-TextCell $ws "E17" "0.43%"

Set-TextCell $ws "D18" "3.450"
Set-TextCell $ws "E18" "10.77%"

Set-TextCell $ws "E19" "0.36%"

Set-TextCell $ws "E20" "-0.58%"

Set-TextCell $ws "D21" "4.827"
Set-TextCell $ws "E21" "-6.84%"

Set-TextCell $ws "E22" "-7.43%"

Set-TextCell $ws "D23" "0.04394"
Set-TextCell $ws "E23" "-0.57%"

Set-TextCell $ws "D24" "0.001236"
Set-TextCell $ws "E24" "0.01%"

Set-TextCell $ws "D25" "0.004180"
Set-TextCell $ws "E25" "-11.27%"

Set-TextCell $ws "E27" "0.03%"

Set-TextCell $ws "D39" "0.02044"
Set-TextCell $ws "E39" "4.84%"

Set-TextCell $ws "D40" "0.05122"
Set-TextCell $ws "E40" "-1.14%"

Set-TextCell $ws "D41" "0.007473"
Set-TextCell $ws "E41" "-1.18%"

Set-TextCell $ws "D42" "0.01002"
Set-TextCell $ws "E42" "-1.43%"

Set-TextCell $ws "D43" "0.1372"
Set-TextCell $ws "E43" "0.08%"

Set-TextCell $ws "D44" "0.002124"
Set-TextCell $ws "E44" "0.98%"

Set-TextCell $ws "D45" "0.009870"
Set-TextCell $ws "E45" "-7.93%"

Set-TextCell $ws "D46" "0.00006318"
Set-TextCell $ws "E46" "-0.19%"

Set-TextCell $ws "E47" "0.03%"

Set-TextCell $ws "D48" "64.82"
Set-TextCell $ws "E48" "-0.61%"

Set-TextCell $ws "E49" "-3.51%"

Set-TextCell $ws "E50" "0.03%"

Set-TextCell $ws "E51" "0.03%"
